$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D for the new 2018-12-31 fiscal year;
# this shifts existing columns D:K to E:L (values, formulas & styles move with them).
$ws.Columns("D").Insert()

# The freshly inserted column D inherits column C's style; copy the
# number formats back from column E (the original D) so D looks right again.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the column width used by the neighbouring data columns.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column with the FY2018 (period ending 2018-12-31) figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 900
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = 47700
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 65900
$ws.Range("D18").Value = -65000
$ws.Range("D20").Value = 1200
$ws.Range("D21").Value = -62900
$ws.Range("D22").Value = 1500
$ws.Range("D23").Value = -65300
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -65300
$ws.Range("D27").Value = -65300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1200
$ws.Range("D33").Value = -65300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -65300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 37400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = "NA"
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 4700
$ws.Range("D46").Value = 42100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 2100
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 44500
$ws.Range("D57").Value = 1100
$ws.Range("D58").Value = 21400
$ws.Range("D59").Value = 12700
$ws.Range("D60").Value = 35100
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 14800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 49900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -280400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -5400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -65300
$ws.Range("D83").Value = 1000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -59200
$ws.Range("D91").Value = -900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 25300
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 700
$ws.Range("D101").Value = -200
$ws.Range("D102").Value = -33300
